$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.07
$ws.Range("H2").Value = 2.77
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2.8
$ws.Range("K2").Value = 1.85
$ws.Range("L2").Value = 4.55
$ws.Range("M2").Value = 1.38
$ws.Range("N2").Value = 2.57
$ws.Range("O2").Value = 2.12
$ws.Range("P2").Value = 1.57
$ws.Range("Q2").Value = 3.45
$ws.Range("R2").Value = 1.22
$ws.Range("S2").Value = 1.52
$ws.Range("T2").Value = 2.22
$ws.Range("U2").Value = 1.8
$ws.Range("V2").Value = 1.82
$ws.Range("W2").Value = 6.3
$ws.Range("AA2").Value = 18.5
$ws.Range("AB2").Value = 30
$ws.Range("AC2").Value = 7
$ws.Range("AD2").Value = 5.5
$ws.Range("AE2").Value = 14
$ws.Range("AF2").Value = 70
$ws.Range("AG2").Value = 600
$ws.Range("AH2").Value = 10
$ws.Range("AI2").Value = 23
$ws.Range("AK2").Value = 70
$ws.Range("AL2").Value = 40
$ws.Range("AM2").Value = 45

# Row 3
$ws.Range("G3").Value = 1.91
$ws.Range("H3").Value = 3.35
$ws.Range("I3").Value = 3.7
$ws.Range("J3").Value = 2.45
$ws.Range("K3").Value = 2.15
$ws.Range("L3").Value = 4.1
$ws.Range("M3").Value = 1.34
$ws.Range("N3").Value = 2.77
$ws.Range("O3").Value = 1.98
$ws.Range("P3").Value = 1.65
$ws.Range("Q3").Value = 3.25
$ws.Range("R3").Value = 1.25
$ws.Range("U3").Value = 1.83
$ws.Range("V3").Value = 1.78
$ws.Range("W3").Value = 6.5
$ws.Range("X3").Value = 8.5
$ws.Range("Y3").Value = 8.5
$ws.Range("Z3").Value = 16
$ws.Range("AA3").Value = 16.5
$ws.Range("AC3").Value = 8.5
$ws.Range("AD3").Value = 6.5
$ws.Range("AE3").Value = 16
$ws.Range("AF3").Value = 80
$ws.Range("AG3").Value = 800
$ws.Range("AH3").Value = 9.75
$ws.Range("AI3").Value = 19.5
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 55
$ws.Range("AL3").Value = 37
$ws.Range("AM3").Value = 45

# Row 4
$ws.Range("G4").Value = 1.21
$ws.Range("H4").Value = 5.7
$ws.Range("J4").Value = 1.6
$ws.Range("K4").Value = 2.62
$ws.Range("P4").Value = 2.22
$ws.Range("Q4").Value = 2.22
$ws.Range("R4").Value = 1.52
$ws.Range("U4").Value = 2.15
$ws.Range("V4").Value = 1.55
$ws.Range("X4").Value = 5.8
$ws.Range("Y4").Value = 9.5
$ws.Range("Z4").Value = 6.7
$ws.Range("AA4").Value = 11
$ws.Range("AB4").Value = 35
$ws.Range("AC4").Value = 14.5
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 29
$ws.Range("AJ4").Value = 40
$ws.Range("AM4").Value = 150

# Row 5
$ws.Range("G5").Value = 2.45
$ws.Range("H5").Value = 2.87
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 3.15
$ws.Range("K5").Value = 1.9
$ws.Range("L5").Value = 3.65
$ws.Range("M5").Value = 1.52
$ws.Range("N5").Value = 2.22
$ws.Range("O5").Value = 2.47
$ws.Range("P5").Value = 1.42
$ws.Range("Q5").Value = 4.25
$ws.Range("R5").Value = 1.14
$ws.Range("S5").Value = 1.52
$ws.Range("T5").Value = 2.2
$ws.Range("U5").Value = 2.07
$ws.Range("V5").Value = 1.6
$ws.Range("W5").Value = 5.9
$ws.Range("X5").Value = 10.25
$ws.Range("Y5").Value = 10.25
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 26
$ws.Range("AB5").Value = 50
$ws.Range("AC5").Value = 6.2
$ws.Range("AE5").Value = 18.5
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 14
$ws.Range("AJ5").Value = 11.5
$ws.Range("AK5").Value = 40
$ws.Range("AL5").Value = 32
$ws.Range("AN5").Value = 1.11
$ws.Range("AO5").Value = 6

# Row 6
$ws.Range("G6").Value = 2.45
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 2.62
$ws.Range("J6").Value = 3
$ws.Range("L6").Value = 3.15
$ws.Range("M6").Value = 1.33
$ws.Range("N6").Value = 2.8
$ws.Range("O6").Value = 1.98
$ws.Range("P6").Value = 1.65
$ws.Range("W6").Value = 7.7
$ws.Range("X6").Value = 11.5
$ws.Range("Y6").Value = 9.75
$ws.Range("Z6").Value = 25
$ws.Range("AA6").Value = 21
$ws.Range("AH6").Value = 8
$ws.Range("AI6").Value = 12.5
$ws.Range("AJ6").Value = 10
$ws.Range("AK6").Value = 29
$ws.Range("AL6").Value = 23
$ws.Range("AM6").Value = 35

# Row 7
$ws.Range("G7").Value = 1.8
$ws.Range("H7").Value = 3.6
$ws.Range("I7").Value = 3.85
$ws.Range("J7").Value = 2.37
$ws.Range("K7").Value = 2.15
$ws.Range("L7").Value = 4.25
$ws.Range("M7").Value = 1.26
$ws.Range("N7").Value = 3.15
$ws.Range("O7").Value = 1.78
$ws.Range("P7").Value = 1.83
$ws.Range("Q7").Value = 2.8
$ws.Range("R7").Value = 1.33
$ws.Range("T7").Value = 2.6
$ws.Range("U7").Value = 1.72
$ws.Range("V7").Value = 1.9
$ws.Range("W7").Value = 7.4
$ws.Range("X7").Value = 8.75
$ws.Range("Y7").Value = 8.25
$ws.Range("Z7").Value = 15
$ws.Range("AA7").Value = 14
$ws.Range("AB7").Value = 25
$ws.Range("AC7").Value = 10.75
$ws.Range("AD7").Value = 7
$ws.Range("AE7").Value = 15
$ws.Range("AH7").Value = 11.5
$ws.Range("AI7").Value = 21
$ws.Range("AJ7").Value = 13
$ws.Range("AK7").Value = 55
$ws.Range("AL7").Value = 35
$ws.Range("AM7").Value = 40

# Row 10
$ws.Range("G10").Value = 2.37
$ws.Range("I10").Value = 3.1
$ws.Range("M10").Value = 1.27
$ws.Range("O10").Value = 2.05
$ws.Range("P10").Value = 1.75
$ws.Range("R10").Value = 1.25
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("U10").Value = 1.8
$ws.Range("V10").Value = 1.91
$ws.Range("X10").Value = 11
$ws.Range("AC10").Value = 9
$ws.Range("AF10").Value = 51
$ws.Range("AG10").Value = 251
$ws.Range("AL10").Value = 26
$ws.Range("AN10").Value = 1.04

# Row 11
$ws.Range("M11").Value = 1.47
$ws.Range("N11").Value = 2.5
$ws.Range("Q11").Value = 4.5
$ws.Range("R11").Value = 1.13
$ws.Range("AN11").Value = 1.08
$ws.Range("AP11").Value = 2
$ws.Range("AQ11").Value = 1.85

# Row 12
$ws.Range("G12").Value = 3.3
$ws.Range("I12").Value = 2.3
$ws.Range("J12").Value = 4
$ws.Range("L12").Value = 3.2
$ws.Range("M12").Value = 1.41
$ws.Range("N12").Value = 2.62
$ws.Range("R12").Value = 1.13
$ws.Range("Z12").Value = 34
$ws.Range("AH12").Value = 6.5
$ws.Range("AI12").Value = 10
$ws.Range("AN12").Value = 1.07

# Row 13
$ws.Range("G13").Value = 2.75
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 2.42
$ws.Range("J13").Value = 3.3
$ws.Range("K13").Value = 2.07
$ws.Range("N13").Value = 3.35
$ws.Range("O13").Value = 1.83
$ws.Range("P13").Value = 1.87
$ws.Range("Q13").Value = 2.95
$ws.Range("W13").Value = 9.75
$ws.Range("X13").Value = 15.5
$ws.Range("Z13").Value = 35
$ws.Range("AA13").Value = 22
$ws.Range("AC13").Value = 7.3
$ws.Range("AD13").Value = 6.2
$ws.Range("AE13").Value = 12
$ws.Range("AH13").Value = 8.5
$ws.Range("AI13").Value = 12.5
$ws.Range("AO13").Value = 7.3

# Row 14
$ws.Range("G14").Value = 2.7
$ws.Range("I14").Value = 2.25
$ws.Range("J14").Value = 3.25
$ws.Range("K14").Value = 2.25
$ws.Range("L14").Value = 2.8
$ws.Range("M14").Value = 1.22
$ws.Range("N14").Value = 3.85
$ws.Range("O14").Value = 1.65
$ws.Range("P14").Value = 2.1
$ws.Range("Q14").Value = 2.6
$ws.Range("R14").Value = 1.44
$ws.Range("S14").Value = 1.33
$ws.Range("T14").Value = 3.05
$ws.Range("U14").Value = 1.57
$ws.Range("V14").Value = 2.25
$ws.Range("W14").Value = 11
$ws.Range("X14").Value = 15.5
$ws.Range("Y14").Value = 10
$ws.Range("Z14").Value = 32
$ws.Range("AB14").Value = 26
$ws.Range("AC14").Value = 8.25
$ws.Range("AD14").Value = 7
$ws.Range("AE14").Value = 12.5
$ws.Range("AF14").Value = 45
$ws.Range("AG14").Value = 300
$ws.Range("AH14").Value = 10
$ws.Range("AI14").Value = 12.5
$ws.Range("AL14").Value = 16.5
$ws.Range("AM14").Value = 23
$ws.Range("AO14").Value = 8.25
